$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename teams in column A (rows 2-4)
$ws.Range("A2").Value = "Team 1"
$ws.Range("A3").Value = "Team 2"
$ws.Range("A4").Value = "Team 3"

# Update the active selection to D6
$ws.Range("D6").Select()
